# Apply updated crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.060.17"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "2.668.29"
$ws.Range("E3").Value = "  +6.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'113.72"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").Value = "'327.58"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "'0.531"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").Value = "'41.43"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").Value = "'20.27"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "3.085.61"
$ws.Range("E15").Value = "  +6.18%  "
$ws.Range("D16").Value = "2.684.68"
$ws.Range("E16").Value = "  +6.79%  "
$ws.Range("D17").Value = "'0.878"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").Value = "49.995.22"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").Value = "'13.34"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'6.82"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").Value = "'2.94"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").Value = "0.0₃0970"
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'72.90"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'279.51"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "'2.60"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").Value = "'27.03"
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D28").Value = "'36.85"
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.95"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.143"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'50.48"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "'19.75"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").Value = "'5.47"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "'0.0803"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "'2.09"
$ws.Range("E37").Value = "  +6.80%  "
$ws.Range("D38").Value = "'4.80"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").Value = "'3.10"
$ws.Range("E39").Value = "  +6.50%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.113"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'22.66"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'123.15"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").Value = "'2.25"
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("D44").Value = "'0.0316"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").Value = "'3.38"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "2.081.63"
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").Value = "'2.25"
$ws.Range("E47").Value = "  +12.39%  "
$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("D49").Value = "'9.19"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").Value = "'5.45"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").Value = "'82.60"
$ws.Range("E51").Value = "  +3.69%  "
